# Update leve-profit data cells across all sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1030.8594
$ws.Range("J17").Value = 1052.9454
$ws.Range("L17").Value = 3158.8362
$ws.Range("N17").Value = -3494.8362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 302062.2
$ws.Range("I61").Value = 225881.19
$ws.Range("K61").Value = 225881.19
$ws.Range("M61").Value = -225669.19

$ws.Range("H74").Value = 122307.19
$ws.Range("I74").Value = 132326.58
$ws.Range("J74").Value = 71542.266
$ws.Range("K74").Value = 132326.58
$ws.Range("L74").Value = 71542.266
$ws.Range("M74").Value = -131452.58
$ws.Range("N74").Value = -73290.266

$ws.Range("H77").Value = 122307.19
$ws.Range("I77").Value = 132326.58
$ws.Range("J77").Value = 71542.266
$ws.Range("K77").Value = 661632.8999999999
$ws.Range("L77").Value = 357711.33
$ws.Range("M77").Value = -657264.8999999999
$ws.Range("N77").Value = -366447.33

$ws.Range("H97").Value = 1007
$ws.Range("I97").Value = 1054.6471
$ws.Range("J97").Value = 737
$ws.Range("K97").Value = 1054.6471
$ws.Range("L97").Value = 737
$ws.Range("M97").Value = -558.6470999999999
$ws.Range("N97").Value = -1729

$ws.Range("H102").Value = 2938.875
$ws.Range("I102").Value = 1700
$ws.Range("K102").Value = 1700
$ws.Range("M102").Value = -78

$ws.Range("H122").Value = 4940.8237
$ws.Range("I122").Value = 4971
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 14913
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -12463
$ws.Range("N122").Value = -19300

$ws.Range("H124").Value = 15429
$ws.Range("J124").Value = 15429
$ws.Range("L124").Value = 15429
$ws.Range("N124").Value = -25249

$ws.Range("H132").Value = 2077.1462
$ws.Range("I132").Value = 1862.6613
$ws.Range("J132").Value = 2742.05
$ws.Range("K132").Value = 5587.9839
$ws.Range("L132").Value = 8226.150000000001
$ws.Range("M132").Value = -3057.9839
$ws.Range("N132").Value = -13286.15

$ws.Range("H136").Value = 302062.2
$ws.Range("I136").Value = 225881.19
$ws.Range("K136").Value = 677643.5700000001
$ws.Range("M136").Value = -675093.5700000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1455.8485
$ws.Range("I20").Value = 1227.5714
$ws.Range("J20").Value = 1855.3334
$ws.Range("K20").Value = 1227.5714
$ws.Range("L20").Value = 1855.3334
$ws.Range("M20").Value = -980.5714
$ws.Range("N20").Value = -2349.3334

$ws.Range("H105").Value = 2011.8276
$ws.Range("I105").Value = 1495.909
$ws.Range("J105").Value = 3633.2856
$ws.Range("K105").Value = 1495.909
$ws.Range("L105").Value = 3633.2856
$ws.Range("M105").Value = 251.0909999999999
$ws.Range("N105").Value = -7127.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2199.8396
$ws.Range("I58").Value = 2260.0967
$ws.Range("J58").Value = 2003.2106
$ws.Range("K58").Value = 2260.0967
$ws.Range("L58").Value = 2003.2106
$ws.Range("M58").Value = -2057.0967
$ws.Range("N58").Value = -2409.2106

$ws.Range("H107").Value = 285.29166
$ws.Range("I107").Value = 136.5
$ws.Range("J107").Value = 731.6667
$ws.Range("K107").Value = 136.5
$ws.Range("L107").Value = 731.6667
$ws.Range("M107").Value = 1783.5
$ws.Range("N107").Value = -4571.6667

$ws.Range("H132").Value = 1331.3
$ws.Range("I132").Value = 755.4865
$ws.Range("J132").Value = 2970.1538
$ws.Range("K132").Value = 2266.4595
$ws.Range("L132").Value = 8910.4614
$ws.Range("M132").Value = 263.5405000000001
$ws.Range("N132").Value = -13970.4614

$ws.Range("H134").Value = 1275.9697
$ws.Range("I134").Value = 778.6
$ws.Range("J134").Value = 2830.25
$ws.Range("K134").Value = 2335.8
$ws.Range("L134").Value = 8490.75
$ws.Range("M134").Value = 199.1999999999998
$ws.Range("N134").Value = -13560.75

$ws.Range("H136").Value = 2199.8396
$ws.Range("I136").Value = 2260.0967
$ws.Range("J136").Value = 2003.2106
$ws.Range("K136").Value = 6780.2901
$ws.Range("L136").Value = 6009.6318
$ws.Range("M136").Value = -4230.2901
$ws.Range("N136").Value = -11109.6318

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 324.85715
$ws.Range("I68").Value = 262.25
$ws.Range("J68").Value = 408.33334
$ws.Range("K68").Value = 786.75
$ws.Range("L68").Value = 1225.00002
$ws.Range("M68").Value = 24.25
$ws.Range("N68").Value = -2847.00002

$ws.Range("H71").Value = 324.85715
$ws.Range("I71").Value = 262.25
$ws.Range("J71").Value = 408.33334
$ws.Range("K71").Value = 2360.25
$ws.Range("L71").Value = 3675.00006
$ws.Range("M71").Value = 1695.75
$ws.Range("N71").Value = -11787.00006

$ws.Range("H92").Value = 510.3
$ws.Range("I92").Value = 433.83334
$ws.Range("J92").Value = 625
$ws.Range("K92").Value = 1301.50002
$ws.Range("L92").Value = 1875
$ws.Range("M92").Value = -53.50001999999995
$ws.Range("N92").Value = -4371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4791.6816
$ws.Range("I107").Value = 8473.583000000001
$ws.Range("J107").Value = 373.4
$ws.Range("K107").Value = 8473.583000000001
$ws.Range("L107").Value = 373.4
$ws.Range("M107").Value = -6553.583000000001
$ws.Range("N107").Value = -4213.4

$ws.Range("H122").Value = 1509.0714
$ws.Range("I122").Value = 1192.7
$ws.Range("K122").Value = 3578.1
$ws.Range("M122").Value = -1128.1

$ws.Range("H126").Value = 2296.1667
$ws.Range("I126").Value = 1458.8572
$ws.Range("J126").Value = 3468.4
$ws.Range("K126").Value = 4376.571599999999
$ws.Range("L126").Value = 10405.2
$ws.Range("M126").Value = -1906.571599999999
$ws.Range("N126").Value = -15345.2

$ws.Range("H132").Value = 2810.157
$ws.Range("I132").Value = 2488.1843
$ws.Range("J132").Value = 3751.3076
$ws.Range("K132").Value = 7464.5529
$ws.Range("L132").Value = 11253.9228
$ws.Range("M132").Value = -4934.5529
$ws.Range("N132").Value = -16313.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 18000
$ws.Range("I88").Value = 18000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 18000
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("M88").Value = -17572

$ws.Range("H91").Value = 18000
$ws.Range("I91").Value = 18000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 18000
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("M91").Value = -16518

$ws.Range("H93").Value = 1446.641
$ws.Range("I93").Value = 1379.25
$ws.Range("J93").Value = 1618.1818
$ws.Range("K93").Value = 1379.25
$ws.Range("L93").Value = 1618.1818
$ws.Range("M93").Value = -131.25
$ws.Range("N93").Value = -4114.1818

$ws.Range("H100").Value = 62505080
$ws.Range("I100").Value = 11840
$ws.Range("J100").Value = 90911096
$ws.Range("K100").Value = 11840
$ws.Range("L100").Value = 90911096
$ws.Range("M100").Value = -11299
$ws.Range("N100").Value = -90912178

$ws.Range("H132").Value = 5136.2856
$ws.Range("I132").Value = 1741.7567
$ws.Range("J132").Value = 11746.685
$ws.Range("K132").Value = 5225.2701
$ws.Range("L132").Value = 35240.055
$ws.Range("M132").Value = -2695.2701
$ws.Range("N132").Value = -40300.055

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3412.5715
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 3722
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3722
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -6468

$ws.Range("H122").Value = 1063.8235
$ws.Range("I122").Value = 959
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 2877
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -427
$ws.Range("N122").Value = -10450

$ws.Range("H132").Value = 1604.26
$ws.Range("I132").Value = 1002.65717
$ws.Range("J132").Value = 3008
$ws.Range("K132").Value = 3007.97151
$ws.Range("L132").Value = 9024
$ws.Range("M132").Value = -477.9715099999999
$ws.Range("N132").Value = -14084
